$d = $word.ActiveDocument

# Namespace/package wrapper used so InsertXML receives a full OOXML package
# fragment (as documented by the runtime's error messages). We target the
# "/word/document.xml" part and supply a <w:body> with the replacement
# <w:p>...</w:p> element(s); InsertXML replaces the exact contents of the
# target Range.
function New-PackageXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Locates the single body paragraph whose text starts with $needle and
# returns its Paragraph object (re-resolved fresh so indices/handles are not
# stale after earlier structural edits).
function Find-ParagraphByText([string]$needle) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph starting with: $needle"
    }
    $s = $rng.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $s -and $p.Range.End -gt $s) {
            return $p
        }
    }
    throw "Could not resolve paragraph index for: $needle"
}

# --- 1) Remove the stray lead-in paragraph before the "Добавление заявки" list ---
$p1 = Find-ParagraphByText("В основном алгоритме будут использоваться следующие процессы:")
$p1.Range.Delete()

# --- 2) Remove the stray lead-in paragraph before the "Начало" list ---
$p2 = Find-ParagraphByText("Для создания блок-схемы в соответствии с ГОСТ 19.701")
$p2.Range.Delete()

# --- 3) Move <w:lastRenderedPageBreak/> from the "3. Детализация..." heading
#        run to the "Алгоритм расчета количества заявок:" run ---
$p3 = Find-ParagraphByText("3. Детализация одной функции")
$p3.Range.InsertXML((New-PackageXml '<w:p w:rsidR="0098469F" w:rsidRPr="0098469F" w:rsidRDefault="0098469F" w:rsidP="0098469F"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="2"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="0098469F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="ru-RU"/></w:rPr><w:t>3. Детализация одной функции (расчет количества заявок и среднего времени ремонта)</w:t></w:r></w:p>'))

$p4 = Find-ParagraphByText("Алгоритм расчета количества заявок:")
$p4.Range.InsertXML((New-PackageXml '<w:p w:rsidR="0098469F" w:rsidRPr="0098469F" w:rsidRDefault="0098469F" w:rsidP="0098469F"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="3"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="0098469F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="ru-RU"/></w:rPr><w:lastRenderedPageBreak/><w:t>Алгоритм расчета количества заявок:</w:t></w:r></w:p>'))

# --- 4) Move <w:lastRenderedPageBreak/> from the "Подсчет количества таких
#        заявок." bullet run to the "Выполнение запроса для получения
#        времени..." bullet run ---
$p5 = Find-ParagraphByText("Подсчет количества таких заявок.")
$p5.Range.InsertXML((New-PackageXml '<w:p w:rsidR="0098469F" w:rsidRPr="0098469F" w:rsidRDefault="0098469F" w:rsidP="0098469F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="0098469F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="ru-RU"/></w:rPr><w:t>Подсчет количества таких заявок.</w:t></w:r></w:p>'))

$p6 = Find-ParagraphByText("Выполнение запроса для получения времени начала и окончания заявок")
$p6.Range.InsertXML((New-PackageXml '<w:p w:rsidR="0098469F" w:rsidRPr="0098469F" w:rsidRDefault="0098469F" w:rsidP="0098469F"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="0098469F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="ru-RU"/></w:rPr><w:lastRenderedPageBreak/><w:t>Выполнение запроса для получения времени начала и окончания заявок со статусом &quot;завершена&quot;.</w:t></w:r></w:p>'))

Write-Output "Edit complete"
